{"js": "async (context) => {\n  const body = context.document.body;\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const secondParagraph = paragraphs.items[1];\n  secondParagraph.insertText(\n    \" Su principio b\u00e1sico de funcionamiento reside en el hecho de que un aumento de la corriente de l\u00ednea provoca una consecuente elevaci\u00f3n de la temperatura debido a las p\u00e9rdidas ocasionadas en el material conductor del dispositivo.\",\n    Word.InsertLocation.end\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$para = $d.Paragraphs(2).Range\n$para.InsertAfter(\" Su principio b\u00e1sico de funcionamiento reside en el hecho de que un aumento de la corriente de l\u00ednea provoca una consecuente elevaci\u00f3n de la temperatura debido a las p\u00e9rdidas ocasionadas en el material conductor del dispositivo.\")\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
